# Restore C10 on the "Rules" sheet back to its earlier value of 1
# (was showing 18 after a later edit; this revision restores 1).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("C10").Value = 1
